$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: (MITRE, 2014) -- B13 becomes italic rich text for the title/source ---
$ws.Range("B13").Value = "MITRE. (2014). CVE - CVE-2014-1532. cve.mitre.org. Retrieved March 18, 2015, from http://cve.mitre.org/cgi-bin/cvename.cgi?name=CVE-2014-1532"
$ws.Range("B13").Characters(16, 19).Font.Italic = $true
$ws.Range("B13").Characters(35, 2).Font.Italic = $false
$ws.Range("B13").Characters(37, 13).Font.Italic = $true
$ws.Range("B13").Characters(50, 92).Font.Italic = $false
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("B13").HorizontalAlignment = 1

# --- Row 14: (SecurityFocus, 2010) ---
$ws.Range("B14").Value = "SecurityFocus. (2010). Adobe Reader and Acrobat U3D Support Remote Code Execution Vulnerability. Retrieved March 11, 2015, from http://www.securityfocus.com/bid/37756/info"
$ws.Range("B14").Characters(24, 72).Font.Italic = $true
$ws.Range("B14").Characters(96, 76).Font.Italic = $false
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").HorizontalAlignment = 1

# --- Row 15: (SecurityTracker, 2010) ---
$ws.Range("B15").Value = "SecurityTracker. (2010). Adobe Acrobat and Adobe Reader Flaws Lets Remote Users Execute Arbitrary Code and Deny Service - SecurityTracker. SecurityTracker. Retrieved March 11, 2015, from http://www.securitytracker.com/id?1023446"
$ws.Range("B15").Characters(26, 112).Font.Italic = $true
$ws.Range("B15").Characters(138, 2).Font.Italic = $false
$ws.Range("B15").Characters(140, 15).Font.Italic = $true
$ws.Range("B15").Characters(155, 74).Font.Italic = $false
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("B15").HorizontalAlignment = 1

# --- Row 16: (MITRE, 2006) ---
$ws.Range("B16").Value = "MITRE. (2006). CVE - CVE-2006-2198. cve.mitre.org. Retrieved March 11, 2015, from http://cve.mitre.org/cgi-bin/cvename.cgi?name=CVE-2006-2198"
$ws.Range("B16").Characters(16, 19).Font.Italic = $true
$ws.Range("B16").Characters(35, 2).Font.Italic = $false
$ws.Range("B16").Characters(37, 13).Font.Italic = $true
$ws.Range("B16").Characters(50, 92).Font.Italic = $false
$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("B16").HorizontalAlignment = 1

# --- Row 17: (Oliver & Snowden, 2015) -- also drop "[Last Week Tonight]" from the author text ---
$ws.Range("B17").Value = "Oliver, J., & Snowden, E. (2015). Last Week Tonight with John Oliver: Edward Snowden on Passwords. Retrieved May 6, 2015, from https://www.youtube.com/watch?v=yzGzB-yYKcc"
$ws.Range("B17").Characters(35, 63).Font.Italic = $true
$ws.Range("B17").Characters(98, 73).Font.Italic = $false
$ws.Range("B17").VerticalAlignment = -4108
$ws.Range("B17").HorizontalAlignment = 1

# --- New row 18: (Goodin, 2014) ---
$ws.Range("A18").Value = "(Goodin, 2014)"
$ws.Range("B18").Value = "Goodin, D. (2014). Stanford’s password policy shuns one-size-fits-all security | Ars Technica. Ars Technica. Retrieved April 30, 2015, from http://arstechnica.com/security/2014/04/25/stanfords-password-policy-shuns-one-size-fits-all-security/"
$ws.Range("B18").Characters(20, 74).Font.Italic = $true
$ws.Range("B18").Characters(94, 2).Font.Italic = $false
$ws.Range("B18").Characters(96, 12).Font.Italic = $true
$ws.Range("B18").Characters(108, 135).Font.Italic = $false
$ws.Range("B18").VerticalAlignment = -4108
$ws.Range("B18").HorizontalAlignment = 1

# --- Selection moves down one row to reflect the newly added entry ---
$ws.Range("B19").Select()
